$wb = $excel.ActiveWorkbook

# The localization status report was regenerated. Two files (1869e430... and
# 4cb333e9...) moved from "Ready for handoff" to "In Translation" for both the
# zh-cn and de-de locales. The c4b19174... file's status is unchanged.

# --- zh-cn sheet (table1 / "Status" column C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de sheet (table2 / "Status" column C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"

# --- Overview sheet: per-locale status summary columns (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
